$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Match formatting/style of the preceding row so the new row looks consistent,
# before putting values in (so no new number-format style gets created).
$ws.Range("A51:I51").Copy() | Out-Null
$ws.Range("A52:I52").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Add the new row of data (row 52), continuing the daily log started in row 2.
$ws.Range("A52").Value = 46009
$ws.Range("B52").Value = 5606
$ws.Range("C52").Value = 4253
$ws.Range("D52").Value = 3882
$ws.Range("E52").Value = 250
$ws.Range("F52").Value = 62
$ws.Range("G52").Value = 49
$ws.Range("H52").Value = 10
$ws.Range("I52").Value = 0

# Update the active selection to the newly added row, as the original file shows.
$ws.Range("A52:I52").Select() | Out-Null
